$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update column G ("K") values for rows 2-8 per the regenerated save data
$ws.Range("G2").Value = 4
$ws.Range("G3").Value = 9
$ws.Range("G4").Value = 2
$ws.Range("G5").Value = 7
$ws.Range("G6").Value = 9
$ws.Range("G7").Value = 4
$ws.Range("G8").Value = 8
